# Fruta / hortaliza, semanal
#
# A new week of price observations (week of serial date 44551) is
# inserted at the top of the "Femacal de La Calera - Naranja" data
# block (rows 465-466), pushing the existing rows 465-580 down to
# 467-582.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 465; everything below
# (old rows 465:580) shifts down to 467:582.
$ws.Rows("465:466").Insert()

# ---- New row 465 --------------------------------------------------
$ws.Cells.Item(465, 1).Value = 3
$ws.Cells.Item(465, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(465, 3).Value = "Coquimbo"
$ws.Cells.Item(465, 4).Value = 44551
$ws.Cells.Item(465, 5).Value = 5
$ws.Cells.Item(465, 6).Value = "Fruta"
$ws.Cells.Item(465, 7).Value = 100102
$ws.Cells.Item(465, 8).Value = "Cítricos"
$ws.Cells.Item(465, 9).Value = 100102005
$ws.Cells.Item(465, 10).Value = "Naranja"
$ws.Cells.Item(465, 11).Value = "Lane Late"
$ws.Cells.Item(465, 12).Value = "Primera"
$ws.Cells.Item(465, 13).Value = 210
$ws.Cells.Item(465, 14).Value = 6000
$ws.Cells.Item(465, 15).Value = 7000
$ws.Cells.Item(465, 16).Value = 6583
$ws.Cells.Item(465, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(465, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(465, 19).Value = 506
$ws.Cells.Item(465, 20).Value = 13

# ---- New row 466 --------------------------------------------------
$ws.Cells.Item(466, 1).Value = 3
$ws.Cells.Item(466, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(466, 3).Value = "Coquimbo"
$ws.Cells.Item(466, 4).Value = 44551
$ws.Cells.Item(466, 5).Value = 5
$ws.Cells.Item(466, 6).Value = "Fruta"
$ws.Cells.Item(466, 7).Value = 100102
$ws.Cells.Item(466, 8).Value = "Cítricos"
$ws.Cells.Item(466, 9).Value = 100102005
$ws.Cells.Item(466, 10).Value = "Naranja"
$ws.Cells.Item(466, 11).Value = "Lane Late"
$ws.Cells.Item(466, 12).Value = "Segunda"
$ws.Cells.Item(466, 13).Value = 120
$ws.Cells.Item(466, 14).Value = 4500
$ws.Cells.Item(466, 15).Value = 5000
$ws.Cells.Item(466, 16).Value = 4792
$ws.Cells.Item(466, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(466, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(466, 19).Value = 369
$ws.Cells.Item(466, 20).Value = 13
